## Server Security Setup Input Variables - update Test1 user-security list
## Mirrors PR 183: Updated Test1 user-security list (related work item #4490)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Security Groups")

# --- Add a "Notes" column (E) with documentation -----------------------
$ws.Range("E1").Value = "Notes"
$ws.Columns.Item(5).ColumnWidth = 39

# Note attached to the Dev / IIS_IUSRS / SPOReportsDev row (row 13)
$ws.Range("E13").Value = "only spore users are 'SPOReportsDev' and 'SPOReports'"

# --- Expand the "Test1" section so it mirrors the "Dev" section --------
# Before: rows 21-23 hold the IIS_IUSRS entries (MammothTest, IconWebTest,
#         Authenticated Users) right after the Administrators entries end
#         at row 20.
# After:  two more Administrators rows are added (NutriconServiceTest,
#         SPOReportsDev) and two more IIS_IUSRS rows are added
#         (NutriconServiceTest, SPOReportsDev) before "Authenticated
#         Users", growing the table from 23 to 27 rows.

# Insert 2 new rows for the Administrators group, right after row 20
$ws.Range("A21:A22").EntireRow.Insert()

$ws.Range("A21").Value = "Test1"
$ws.Range("B21").Value = "Administrators"
$ws.Range("C21").Value = "WFM\NutriconServiceTest"
$ws.Range("D21").Value = "Web"

$ws.Range("A22").Value = "Test1"
$ws.Range("B22").Value = "Administrators"
$ws.Range("C22").Value = "WFM\SPOReportsDev"
$ws.Range("D22").Value = "Web"

# The old IIS_IUSRS rows (MammothTest, IconWebTest, Authenticated Users)
# have shifted down to rows 23-25. Insert 2 new rows before the
# "Authenticated Users" row (now row 25) to add the matching IIS_IUSRS
# entries.
$ws.Range("A25:A26").EntireRow.Insert()

$ws.Range("A25").Value = "Test1"
$ws.Range("B25").Value = "IIS_IUSRS"
$ws.Range("C25").Value = "WFM\NutriconServiceTest"
$ws.Range("D25").Value = "Web"

$ws.Range("A26").Value = "Test1"
$ws.Range("B26").Value = "IIS_IUSRS"
$ws.Range("C26").Value = "WFM\SPOReportsDev"
$ws.Range("D26").Value = "Web"

# --- Refresh the AutoFilter so it covers the full new range ------------
$ws.AutoFilterMode = $false
$ws.Range("A1:E27").AutoFilter()

# Keep the workbook-level filter-database defined name in sync
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Security Groups!_FilterDatabase") {
        $n.RefersTo = "='Security Groups'!`$A`$1:`$E`$27"
    }
}

# --- Match the saved selection state ------------------------------------
$ws.Activate()
$ws.Range("A28").Select()
